# Insert a new weekly data row for "Perejil" (Feria Lagunitas de Puerto Montt)
# at worksheet row 352, pushing the existing rows 352-376 down to 353-377.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row, shifting rows 352..376 down to 353..377.
$ws.Rows.Item(352).Insert()

# Populate the newly inserted row 352 with the new record.
$ws.Range("A352").Value = 4
$ws.Range("B352").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C352").Value = "Los Lagos"
$ws.Range("D352").Value = 45021
$ws.Range("E352").Value = 10
$ws.Range("F352").Value = 100112044
$ws.Range("G352").Value = "Perejil"
$ws.Range("H352").Value = "Sin especificar"
$ws.Range("I352").Value = "Primera"
$ws.Range("J352").Value = 20
$ws.Range("K352").Value = 6000
$ws.Range("L352").Value = 6000
$ws.Range("M352").Value = 6000
$ws.Range("N352").Value = '$/docena de atados (2 kilos)'
$ws.Range("O352").Value = "Región de La Araucanía"
$ws.Range("P352").Value = 3000
$ws.Range("Q352").Value = 2
$ws.Range("R352").Value = "Hortaliza"

# Keep the style of the date cell (D) consistent with the rest of the column.
$ws.Range("D352").NumberFormat = $ws.Range("D353").NumberFormat
